$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2..4 hold species-observation records; this commit rotates the three
# data rows (2,3,4) up by one: new row2 = old row3, new row3 = old row4,
# new row4 = old row2. Values are written directly (captured from the
# original cells) rather than read back through Range.Value, and only the
# cells whose content actually changes are touched.

# -- Row 2 --
$ws.Range("A2").Value = 80000956
$ws.Range("B2").Value = 57133
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 100041
$ws.Range("F2").Value = "Hasselsnok"
$ws.Range("G2").Value = "Coronella austriaca"
$ws.Range("H2").Value = "Laurenti, 1768"
$ws.Range("I2").Value = "'1"
$ws.Range("J2").Value = "ex."
$ws.Range("N2").Value = "observerad"
$ws.Range("P2").Value = "Ekbacken, Ög"
$ws.Range("Q2").Value = 574245.9331973131
$ws.Range("R2").Value = 6505393.860064601
$ws.Range("S2").Value = 5
$ws.Range("Y2").Value = "'2015-07-13"
$ws.Range("AA2").Value = "'2015-07-13"
$ws.Range("AC2").Value = "observerad"
$ws.Range("AI2").Value = "rasbrant"
$ws.Range("AW2").Value = "Elin Håkansson"
$ws.Range("AX2").Value = "Mikael Hagström"
$ws.Range("AY2").Value = "Ostlänken Norrköpings kommun (OLP2)"
# -- Row 3 --
$ws.Range("A3").Value = 97784188
$ws.Range("B3").Value = 78098
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6453
$ws.Range("F3").Value = "Vedskivlav"
$ws.Range("G3").Value = "Hertelidea botryosa"
$ws.Range("H3").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("P3").Value = "Korpklint, Ög"
$ws.Range("Q3").Value = 574267.5050753297
$ws.Range("R3").Value = 6505328.227546699
$ws.Range("S3").Value = 15
$ws.Range("Y3").Value = "'2021-12-30"
$ws.Range("AA3").Value = "'2021-12-30"
$ws.Range("AC3").Value = ""
$ws.Range("AI3").Value = ""
$ws.Range("AW3").Value = "Marika Sjödin"
$ws.Range("AX3").Value = "Marika Sjödin, Eva Siljeholm"
$ws.Range("AY3").Value = ""
# -- Row 4 --
$ws.Range("A4").Value = 98350752
$ws.Range("B4").Value = 93132
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 2671
$ws.Range("F4").Value = "Fällmossa"
$ws.Range("G4").Value = "Antitrichia curtipendula"
$ws.Range("H4").Value = "(Hedw.) Brid."
$ws.Range("Q4").Value = 573934.5192830344
$ws.Range("R4").Value = 6505660.931744166
$ws.Range("Y4").Value = "'2022-01-26"
$ws.Range("AA4").Value = "'2022-01-26"
